$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix data bug: correct covid_deaths counts and realign date/agegrp rows 833-950 ---

# Row 833 (CHANGED)
$ws.Cells.Item(833, 1).Value = 44126
$ws.Cells.Item(833, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(833, 2).Value = "70-79"
$ws.Cells.Item(833, 3).Value = 9

# Row 842 (CHANGED)
$ws.Cells.Item(842, 1).Value = 44128
$ws.Cells.Item(842, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(842, 2).Value = "70-79"
$ws.Cells.Item(842, 3).Value = 8

# Row 899 (CHANGED)
$ws.Cells.Item(899, 1).Value = 44140
$ws.Cells.Item(899, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(899, 2).Value = "70-79"
$ws.Cells.Item(899, 3).Value = 11

# Row 904 (CHANGED)
$ws.Cells.Item(904, 1).Value = 44141
$ws.Cells.Item(904, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(904, 2).Value = "70-79"
$ws.Cells.Item(904, 3).Value = 14

# Row 926 (CHANGED)
$ws.Cells.Item(926, 1).Value = 44146
$ws.Cells.Item(926, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(926, 2).Value = "80+"
$ws.Cells.Item(926, 3).Value = 28

# Row 927 (CHANGED)
$ws.Cells.Item(927, 1).Value = 44147
$ws.Cells.Item(927, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(927, 2).Value = "40-49"
$ws.Cells.Item(927, 3).Value = 1

# Row 928 (CHANGED)
$ws.Cells.Item(928, 1).Value = 44147
$ws.Cells.Item(928, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(928, 2).Value = "50-59"
$ws.Cells.Item(928, 3).Value = 1

# Row 929 (CHANGED)
$ws.Cells.Item(929, 1).Value = 44147
$ws.Cells.Item(929, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(929, 2).Value = "60-69"
$ws.Cells.Item(929, 3).Value = 3

# Row 930 (CHANGED)
$ws.Cells.Item(930, 1).Value = 44147
$ws.Cells.Item(930, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(930, 2).Value = "70-79"
$ws.Cells.Item(930, 3).Value = 9

# Row 931 (CHANGED)
$ws.Cells.Item(931, 1).Value = 44147
$ws.Cells.Item(931, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(931, 2).Value = "80+"
$ws.Cells.Item(931, 3).Value = 14

# Row 932 (CHANGED)
$ws.Cells.Item(932, 1).Value = 44148
$ws.Cells.Item(932, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(932, 2).Value = "40-49"
$ws.Cells.Item(932, 3).Value = 1

# Row 933 (CHANGED)
$ws.Cells.Item(933, 1).Value = 44148
$ws.Cells.Item(933, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(933, 2).Value = "50-59"
$ws.Cells.Item(933, 3).Value = 1

# Row 934 (CHANGED)
$ws.Cells.Item(934, 1).Value = 44148
$ws.Cells.Item(934, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(934, 2).Value = "60-69"
$ws.Cells.Item(934, 3).Value = 5

# Row 935 (CHANGED)
$ws.Cells.Item(935, 1).Value = 44148
$ws.Cells.Item(935, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(935, 2).Value = "70-79"
$ws.Cells.Item(935, 3).Value = 7

# Row 936 (CHANGED)
$ws.Cells.Item(936, 1).Value = 44148
$ws.Cells.Item(936, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(936, 2).Value = "80+"
$ws.Cells.Item(936, 3).Value = 15

# Row 937 (CHANGED)
$ws.Cells.Item(937, 1).Value = 44149
$ws.Cells.Item(937, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(937, 2).Value = "40-49"
$ws.Cells.Item(937, 3).Value = 1

# Row 938 (CHANGED)
$ws.Cells.Item(938, 1).Value = 44149
$ws.Cells.Item(938, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(938, 2).Value = "50-59"
$ws.Cells.Item(938, 3).Value = 3

# Row 939 (CHANGED)
$ws.Cells.Item(939, 1).Value = 44149
$ws.Cells.Item(939, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(939, 2).Value = "60-69"
$ws.Cells.Item(939, 3).Value = 2

# Row 940 (CHANGED)
$ws.Cells.Item(940, 1).Value = 44149
$ws.Cells.Item(940, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(940, 2).Value = "70-79"
$ws.Cells.Item(940, 3).Value = 9

# Row 941 (CHANGED)
$ws.Cells.Item(941, 1).Value = 44149
$ws.Cells.Item(941, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(941, 2).Value = "80+"
$ws.Cells.Item(941, 3).Value = 20

# Row 942 (CHANGED)
$ws.Cells.Item(942, 1).Value = 44150
$ws.Cells.Item(942, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(942, 2).Value = "40-49"
$ws.Cells.Item(942, 3).Value = 1

# Row 943 (CHANGED)
$ws.Cells.Item(943, 1).Value = 44150
$ws.Cells.Item(943, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(943, 2).Value = "50-59"
$ws.Cells.Item(943, 3).Value = 2

# Row 944 (NEW)
$ws.Cells.Item(944, 1).Value = 44150
$ws.Cells.Item(944, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(944, 2).Value = "60-69"
$ws.Cells.Item(944, 3).Value = 9

# Row 945 (NEW)
$ws.Cells.Item(945, 1).Value = 44150
$ws.Cells.Item(945, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(945, 2).Value = "70-79"
$ws.Cells.Item(945, 3).Value = 12

# Row 946 (NEW)
$ws.Cells.Item(946, 1).Value = 44150
$ws.Cells.Item(946, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(946, 2).Value = "80+"
$ws.Cells.Item(946, 3).Value = 21

# Row 947 (NEW)
$ws.Cells.Item(947, 1).Value = 44151
$ws.Cells.Item(947, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(947, 2).Value = "30-39"
$ws.Cells.Item(947, 3).Value = 1

# Row 948 (NEW)
$ws.Cells.Item(948, 1).Value = 44151
$ws.Cells.Item(948, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(948, 2).Value = "60-69"
$ws.Cells.Item(948, 3).Value = 1

# Row 949 (NEW)
$ws.Cells.Item(949, 1).Value = 44151
$ws.Cells.Item(949, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(949, 2).Value = "70-79"
$ws.Cells.Item(949, 3).Value = 7

# Row 950 (NEW)
$ws.Cells.Item(950, 1).Value = 44151
$ws.Cells.Item(950, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(950, 2).Value = "80+"
$ws.Cells.Item(950, 3).Value = 13
